$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 entirely (the MuSCs -> Resolving-Mac row)
$ws.Rows.Item(5).Delete()

# Update recalculated TPM-based values in rows 2-4 (columns M through T)
$ws.Range("M2").Value = 0.5447316666666667
$ws.Range("N2").Value = 1.634195
$ws.Range("O2").Value = 0.1484165462704666
$ws.Range("P2").Value = 0.1484165462704666
$ws.Range("Q2").Value = 0.05992847273111111
$ws.Range("R2").Value = 0.53935625458
$ws.Range("S2").Value = 0.1484165462704666
$ws.Range("T2").Value = 0.1484165462704666

$ws.Range("O3").Value = 0.6462308875194944
$ws.Range("P3").Value = 0.6462308875194943
$ws.Range("S3").Value = 0.6462308875194944
$ws.Range("T3").Value = 0.6462308875194943

$ws.Range("M4").Value = 0.7537033333333333
$ws.Range("N4").Value = 2.26111
$ws.Range("O4").Value = 0.205352566210039
$ws.Range("P4").Value = 0.205352566210039
$ws.Range("Q4").Value = 0.08291842098222221
$ws.Range("R4").Value = 0.7462657888399999
$ws.Range("S4").Value = 0.205352566210039
$ws.Range("T4").Value = 0.205352566210039
